$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 142 (shifts existing rows 142-177 down to 143-178,
# and extends the used range from A1:R177 to A1:R178).
$ws.Rows(142).Insert()

# Populate the newly inserted row 142 with the new weekly price-report record.
$ws.Range("A142").Value = 11
$ws.Range("B142").Value = 'Vega Monumental Concepción'
$ws.Range("C142").Value = 'Bíobío'
$ws.Range("D142").Value = 44943
$ws.Range("E142").Value = 8
$ws.Range("F142").Value = 100112043
$ws.Range("G142").Value = 'Pepino ensalada'
$ws.Range("H142").Value = 'Sin especificar'
$ws.Range("I142").Value = 'Primera'
$ws.Range("J142").Value = 350
$ws.Range("K142").Value = 10000
$ws.Range("L142").Value = 11000
$ws.Range("M142").Value = 10429
$ws.Range("N142").Value = '$/caja 60 unidades'
$ws.Range("O142").Value = 'Región Metropolitana'
$ws.Range("P142").Value = 174
$ws.Range("Q142").Value = 60
$ws.Range("R142").Value = 'Hortaliza'
